# Apply the "Created Separated classes to build POJO for user and batch
# module" test-data updates to TestData.xlsx:
#   - Batch sheet: add a new "CreateBatchWithValidDataFP" scenario row right
#     after the header, and tweak the BatchName/ProgramId of the existing
#     "CreateBatchWithValidData" scenario.
#   - User sheet: give the "CreateUserWithValidDataR01" scenario a fresh
#     phone number / login email.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Batch sheet
# ---------------------------------------------------------------------
$batch = $wb.Worksheets.Item("Batch")

# Insert a new blank row at the very top (row 1); this shifts the header
# (and every data row below it) down by one.
$batch.Rows(1).Insert()

# Insert a second new blank row right after the header (now row 2); this
# shifts the former "CreateBatchWithValidData" row (and everything below
# it) down once more, opening up row 3 for the brand new scenario.
$batch.Rows(3).Insert()

# Fill the newly inserted row 3 with the new "CreateBatchWithValidDataFP"
# scenario.
$batch.Range("A3").Value = "CreateBatchWithValidDataFP"
$batch.Range("B3").Value = "Automation11221"
$batch.Range("C3").Value = "UnitTesting11"
$batch.Range("D3").Value = 2
$batch.Range("E3").Value = "Active"
$batch.Range("F3").Value = 39
$batch.Range("G3").Value = 201
$batch.Range("I3").Value = "application/json"
$batch.Range("I3").WrapText = $true

# Update the existing "CreateBatchWithValidData" scenario, now on row 4.
$batch.Range("C4").Value = "Python11112"
$batch.Range("F4").Value = 40

# Record the new selection on the Batch sheet without leaving it as the
# active tab (match the source workbook, which stays on "User").
$batch.Activate()
$batch.Range("C4").Select()

# ---------------------------------------------------------------------
# User sheet
# ---------------------------------------------------------------------
$user = $wb.Worksheets.Item("User")

$user.Range("J2").Value = "+91 1666698881"
$user.Range("P2").Value = "NinjaA34@gmail.com"

$user.Activate()
$user.Range("J2").Select()
